# The "COMPETENCES TECHNIQUES" section lists 7 skill-category lines,
# one per paragraph, all sharing identical paragraph formatting
# (<w:spacing w:line="240" w:lineRule="auto" w:before="0" w:after="0"/>).
# The edit simply reorders these 7 lines; nothing about their formatting
# changes. Rather than physically moving paragraphs/runs around, we locate
# the section heading and rewrite the text of the 7 paragraphs that follow
# it, in place, to match the new order:
#
# Old order:
#   Web : client
#   Langages : r, python, matlab, c, c++
#   Bases de données : SQL, MongoDB, Neo4j, Redis
#   Autres : marketing
#   Visualisation : tableau
#   ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn
#   MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit
#
# New order:
#   Langages : r, python, matlab, c, c++
#   Visualisation : tableau
#   MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit
#   Web : client
#   Autres : marketing
#   ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn
#   Bases de données : SQL, MongoDB, Neo4j, Redis

$d = $word.ActiveDocument

# Locate the "COMPETENCES TECHNIQUES" heading paragraph so the 7 lines
# below it are addressed relative to it (robust to anything earlier in
# the document shifting paragraph indices around).
$headingIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*COMPETENCES TECHNIQUES*") {
        $headingIndex = $i
        break
    }
}

$newLines = @(
    "Langages : r, python, matlab, c, c++",
    "Visualisation : tableau",
    "MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit",
    "Web : client",
    "Autres : marketing",
    "ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn",
    "Bases de données : SQL, MongoDB, Neo4j, Redis"
)

for ($j = 0; $j -lt $newLines.Length; $j++) {
    $paraIndex = $headingIndex + 1 + $j
    $d.Paragraphs($paraIndex).Range.Text = $newLines[$j]
}

Write-Host "Reordered COMPETENCES TECHNIQUES section (heading at paragraph $headingIndex)"
